# Auto-generated edit script for cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-ForcedTextCell($ws, $addr, $val) {
    # Ensure numeric-looking strings (e.g. '211.89', '1.00') are stored as
    # literal text (preserving trailing zeros / exact formatting) rather than
    # being auto-converted to a number by Excel.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '26.750.49'
Set-TextCell $ws 'E2' '  +0.25%  '

# Row 3
Set-TextCell $ws 'D3' '1.603.94'
Set-TextCell $ws 'E3' '  +0.36%  '

# Row 4
Set-TextCell $ws 'E4' '  +0.13%  '

# Row 5
Set-ForcedTextCell $ws 'D5' '211.89'
Set-TextCell $ws 'E5' '  +0.18%  '

# Row 6
Set-TextCell $ws 'E6' '  +0.22%  '

# Row 7
Set-TextCell $ws 'E7' '  +0.19%  '

# Row 8
Set-TextCell $ws 'E8' '  +0.15%  '

# Row 9
Set-TextCell $ws 'E9' '  +0.17%  '

# Row 10
Set-ForcedTextCell $ws 'D10' '19.62'
Set-TextCell $ws 'E10' '  +0.75%  '

# Row 11
Set-ForcedTextCell $ws 'D11' '0.0849'
Set-TextCell $ws 'E11' '  +0.77%  '

# Row 12
Set-TextCell $ws 'D12' '1.829.98'
Set-TextCell $ws 'E12' '  +0.40%  '

# Row 13
Set-TextCell $ws 'D13' '1.601.50'
Set-TextCell $ws 'E13' '  +0.23%  '

# Row 14
Set-TextCell $ws 'E14' '  +1.11%  '

# Row 16
Set-ForcedTextCell $ws 'D16' '65.13'
Set-TextCell $ws 'E16' '  +0.02%  '

# Row 17
Set-TextCell $ws 'D17' '0.0₃0741'
Set-TextCell $ws 'E17' '  -1.39%  '

# Row 18
Set-TextCell $ws 'B18' 'Dai'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-ForcedTextCell $ws 'D18' '1.01'
Set-TextCell $ws 'E18' '  +0.13%  '

# Row 19
Set-ForcedTextCell $ws 'D19' '209.27'
Set-TextCell $ws 'E19' '  -0.30%  '

# Row 20
Set-TextCell $ws 'B20' 'Chainlink'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-ForcedTextCell $ws 'D20' '7.19'
Set-TextCell $ws 'E20' '  +2.16%  '

# Row 21
Set-TextCell $ws 'E21' '  +0.53%  '

# Row 22
Set-TextCell $ws 'E22' '  -4.49%  '

# Row 23
Set-ForcedTextCell $ws 'D23' '9.05'
Set-TextCell $ws 'E23' '  +0.93%  '

# Row 24
Set-ForcedTextCell $ws 'D24' '143.55'
Set-TextCell $ws 'E24' '  +0.15%  '

# Row 25
Set-TextCell $ws 'E25' '  +0.02%  '

# Row 26
Set-ForcedTextCell $ws 'D26' '7.13'
Set-TextCell $ws 'E26' '  +0.31%  '

# Row 27
Set-TextCell $ws 'E27' '  -0.14%  '

# Row 29
Set-ForcedTextCell $ws 'D29' '0.0508'
Set-TextCell $ws 'E29' '  -1.54%  '

# Row 30
Set-TextCell $ws 'E30' '  +0.20%  '

# Row 31
Set-ForcedTextCell $ws 'D31' '3.28'
Set-TextCell $ws 'E31' '  +0.69%  '

# Row 32
Set-TextCell $ws 'E32' '  +0.24%  '

# Row 33
Set-TextCell $ws 'D33' '1.284.07'
Set-TextCell $ws 'E33' '  -0.31%  '

# Row 34
Set-TextCell $ws 'E34' '  +1.54%  '

# Row 35
Set-ForcedTextCell $ws 'D35' '1.23'
Set-TextCell $ws 'E35' '  +16.10%  '

# Row 36
Set-TextCell $ws 'E36' '  +0.29%  '

# Row 37
Set-TextCell $ws 'E37' '  -4.66%  '

# Row 38
Set-TextCell $ws 'E38' '  -1.03%  '

# Row 39
Set-ForcedTextCell $ws 'D39' '0.826'
Set-TextCell $ws 'E39' '  +0.04%  '

# Row 40
Set-ForcedTextCell $ws 'D40' '5.47'
Set-TextCell $ws 'E40' '  +0.49%  '

# Row 41
Set-TextCell $ws 'E41' '  +0.29%  '

# Row 42
Set-ForcedTextCell $ws 'D42' '0.779'
Set-TextCell $ws 'E42' '  -0.37%  '

# Row 43
Set-ForcedTextCell $ws 'D43' '62.70'
Set-TextCell $ws 'E43' '  -0.85%  '

# Row 44
Set-TextCell $ws 'D44' '1.741.73'
Set-TextCell $ws 'E44' '  +0.39%  '

# Row 45
Set-ForcedTextCell $ws 'D45' '90.37'
Set-TextCell $ws 'E45' '  -0.59%  '

# Row 46
Set-TextCell $ws 'E46' '  -0.02%  '

# Row 47
Set-TextCell $ws 'E47' '  +1.91%  '

# Row 48
Set-TextCell $ws 'B48' 'BabyDogeCoin'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 'D48' '0.0₆0103'
Set-TextCell $ws 'E48' '  -2.82%  '

# Row 49
Set-TextCell $ws 'B49' 'Cronos'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-ForcedTextCell $ws 'D49' '0.0512'
Set-TextCell $ws 'E49' '  +0.69%  '

# Row 50
Set-TextCell $ws 'B50' 'EnergySwap'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-ForcedTextCell $ws 'D50' '7.57'
Set-TextCell $ws 'E50' '  +3.15%  '

# Row 51
Set-TextCell $ws 'B51' 'USDD'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-ForcedTextCell $ws 'D51' '1.00'
Set-TextCell $ws 'E51' '  +0.05%  '
